{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same textual changes described by the unified diff:\n//  1) \"O(log(N))\" -> \"O(log(N-1))\" in the \"mayor_menor\" paragraph, and the\n//     trailing explanation sentence is reworded.\n//  2) \"O(log(N))\" -> \"O(log(ra\u00edz cuadrada de N))\" in the \"buscar\" paragraph,\n//     and its trailing explanation sentence is reworded / extended.\n\n// --- Paragraph 1: \"mayor_menor\" explanation -------------------------------\n\n// 1a) \"N\" -> \"N-1\" inside \"O(log(N))\". Scope the search narrowly (match\n// case, whole word) so only the intended \"N\" is touched, not any other\n// letter N that might appear elsewhere in the paragraph text.\nconst nResults = context.document.body.search(\"log(N)\", { matchCase: true });\nnResults.load(\"text\");\nawait context.sync();\n\nif (nResults.items.length > 0) {\n  const hit = nResults.items[0];\n  hit.insertText(\"log(N-1)\", \"Replace\");\n  await context.sync();\n}\n\n// 1b) Reword the trailing explanation sentence.\nconst oldTail1 =\n  \", ya que recorrer\u00edamos la lista una vez por cada elemento restando a su vez un elemento a la lista y por muchos datos k metamos no se prolonga en el tiempo.\";\nconst newTail1 =\n  \", ya que recorrer\u00edamos la lista una vez por cada elemento, pero restando a su vez\u2026 un elemento a la lista.\";\n\nconst tail1Results = context.document.body.search(oldTail1, { matchCase: true });\ntail1Results.load(\"text\");\nawait context.sync();\n\nif (tail1Results.items.length > 0) {\n  tail1Results.items[0].insertText(newTail1, \"Replace\");\n  await context.sync();\n}\n\n// --- Paragraph 2: \"buscar\" explanation -------------------------------------\n\n// 2a) \"N\" -> \"ra\u00edz cuadrada de N\" inside the second \"O(log(N))\".\nconst n2Results = context.document.body.search(\"log(N)\", { matchCase: true });\nn2Results.load(\"text\");\nawait context.sync();\n\nif (n2Results.items.length > 0) {\n  const hit = n2Results.items[0];\n  hit.insertText(\"log(ra\u00edz cuadrada de N)\", \"Replace\");\n  await context.sync();\n}\n\n// 2b) Reword / extend the trailing explanation sentence.\nconst oldTail2 =\n  \"ya que recorre un solo bucle y por mas datos que metamos no se prolonga en el tiempo.\";\nconst newTail2 =\n  \"ya que recorre un solo bucle y por m\u00e1s datos que metamos no se prolonga en el tiempo ya que va reduciendo a la mitad de la mitad, el trabajo de busqueda .\";\n\nconst tail2Results = context.document.body.search(oldTail2, { matchCase: true });\ntail2Results.load(\"text\");\nawait context.sync();\n\nif (tail2Results.items.length > 0) {\n  tail2Results.items[0].insertText(newTail2, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same textual changes described by the unified diff:\n#  1) \"O(log(N))\" -> \"O(log(N-1))\" in the \"mayor_menor\" paragraph, and the\n#     trailing explanation sentence is reworded.\n#  2) \"O(log(N))\" -> \"O(log(ra\u00edz cuadrada de N))\" in the \"buscar\" paragraph,\n#     and its trailing explanation sentence is reworded / extended.\n\n$d = $word.ActiveDocument\n\n# wdReplace constants (literal values; not relying on any pre-seeded\n# $wd* globals): wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-InRange($range, $findText, $replaceText) {\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    return $range.Find.Execute(\n        $findText,      # FindText\n        $false,         # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        $wdFindContinue,# Wrap\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        $wdReplaceAll   # Replace\n    )\n}\n\n# --- Paragraph 1: \"mayor_menor\" explanation --------------------------------\n# Locate the paragraph containing \"mayor_menor\" explicitly rather than\n# assuming a fixed paragraph index, so the script is resilient to minor\n# structural differences.\n$para1 = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*mayor_menor*\") {\n        $para1 = $p\n        break\n    }\n}\n\nif ($para1 -ne $null) {\n    # 1a) \"N\" -> \"N-1\" inside \"O(log(N))\".\n    Replace-InRange $para1.Range \"log(N)\" \"log(N-1)\"\n\n    # 1b) Reword the trailing explanation sentence.\n    $oldTail1 = \", ya que recorrer\u00edamos la lista una vez por cada elemento restando a su vez un elemento a la lista y por muchos datos k metamos no se prolonga en el tiempo.\"\n    $newTail1 = \", ya que recorrer\u00edamos la lista una vez por cada elemento, pero restando a su vez\u2026 un elemento a la lista.\"\n    Replace-InRange $para1.Range $oldTail1 $newTail1\n}\n\n# --- Paragraph 2: \"buscar\" explanation --------------------------------------\n$para2 = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*La funci\u00f3n buscar*\") {\n        $para2 = $p\n        break\n    }\n}\n\nif ($para2 -ne $null) {\n    # 2a) \"N\" -> \"ra\u00edz cuadrada de N\" inside the second \"O(log(N))\".\n    Replace-InRange $para2.Range \"log(N)\" \"log(ra\u00edz cuadrada de N)\"\n\n    # 2b) Reword / extend the trailing explanation sentence.\n    $oldTail2 = \"ya que recorre un solo bucle y por mas datos que metamos no se prolonga en el tiempo.\"\n    $newTail2 = \"ya que recorre un solo bucle y por m\u00e1s datos que metamos no se prolonga en el tiempo ya que va reduciendo a la mitad de la mitad, el trabajo de busqueda .\"\n    Replace-InRange $para2.Range $oldTail2 $newTail2\n}\n\nWrite-Output \"done\"\n"}
